# Edit script: updates Eve LUNG_CampA_timetable.xlsx ("Full Timetable" sheet)
# to reflect the refactored daily-schedule generation (dictionary based),
# the H2/H9 student-name swap and the corrected evening / Saturday-morning
# activity assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Unmerge all the cell ranges (rows 7-14, columns B-E) that need to be
#    reshaped so that values can be written to every individual cell first.
# ---------------------------------------------------------------------------
$ws.Range("B7:B10").UnMerge()
$ws.Range("B11:B14").UnMerge()
$ws.Range("C7:C14").UnMerge()
$ws.Range("D7:D14").UnMerge()
$ws.Range("E7:E10").UnMerge()
$ws.Range("E11:E14").UnMerge()

# ---------------------------------------------------------------------------
# 2. Row 7 (11:00 slot) -- swap "Free Time" / "Practice" for column C and
#    assign the private lesson with Liya HUANG to column D.
# ---------------------------------------------------------------------------
$ws.Range("C7").Value = "Practice `n(Harp practice room)"
$ws.Range("D7").Value = "Private Lesson with Liya HUANG `n(Room 242)"

# ---------------------------------------------------------------------------
# 3. Row 11 (12:00 slot) -- the private lesson with Sivan MEGAN moves from
#    column B to column C, column D now gets "Practice", column B & E become
#    blank (they're absorbed into the B7:B14 / E7:E14 merges).
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = "Private Lesson with Sivan MEGAN `n(Room 245)"
$ws.Range("D11").Value = "Practice `n(Harp practice room)"
$ws.Range("E11").Value = ""

# ---------------------------------------------------------------------------
# 4. Re-merge the rows 7-14 ranges according to the new layout. Re-apply a
#    plain thin border around each merged block so the cells keep using the
#    same "thin box" style as the rest of the sheet (UnMerge/Merge otherwise
#    leaves the per-row border variants behind).
# ---------------------------------------------------------------------------
$newRanges7to14 = @("B7:B14", "C7:C10", "C11:C14", "D7:D10", "D11:D14", "E7:E14")
foreach ($r in $newRanges7to14) {
    $ws.Range($r).Merge()
    $ws.Range($r).Borders.LineStyle = 1
}

# ---------------------------------------------------------------------------
# 5. Row 20 (14:15 slot) -- rotate the Gwyneth WENTINK private lesson /
#    Practice / Free Time values, and give column E the Sivan MEGAN lesson.
# ---------------------------------------------------------------------------
$ws.Range("B20").Value = "Private Lesson with Gwyneth WENTINK `n(Room 236)"
$ws.Range("C20").Value = "Practice `n(Harp practice room)"
$ws.Range("D20").Value = "Practice `n(Harp practice room)"
$ws.Range("E20").Value = "Private Lesson with Sivan MEGAN `n(Room 245)"

# ---------------------------------------------------------------------------
# 6. Row 24 (15:15 slot) -- "Ensemble (Room 245)" becomes "Acting Class
#    (Room G13)" for B-E, while F24 becomes "Group Activity".
# ---------------------------------------------------------------------------
$ws.Range("B24").Value = "Acting Class `n(Room G13)"
$ws.Range("C24").Value = "Acting Class `n(Room G13)"
$ws.Range("D24").Value = "Acting Class `n(Room G13)"
$ws.Range("E24").Value = "Acting Class `n(Room G13)"
$ws.Range("F24").Value = "Group Activity `n(Room Group Activity)"

# ---------------------------------------------------------------------------
# 7. Row 28 (16:15 slot) -- "Acting Class (Room G13)" becomes "Ensemble
#    (Room 236)" for B-E, while F28 becomes "Break".
# ---------------------------------------------------------------------------
$ws.Range("B28").Value = "Ensemble `n(Room 236)"
$ws.Range("C28").Value = "Ensemble `n(Room 236)"
$ws.Range("D28").Value = "Ensemble `n(Room 236)"
$ws.Range("E28").Value = "Ensemble `n(Room 236)"
$ws.Range("F28").Value = "Break"
